# Week 17 data logging + tiebreak-method fix for the Buccaneers 2021 Team Data workbook.
# Appends this week's per-drive/per-play samples to the YDS and ST running logs,
# and bumps the Home/Road season totals on OFF, DEF, ST, TURNS and PEN to include
# the new week's numbers.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# YDS sheet: append Week 17 rushing (R) / passing (P) play-by-play yardage
# samples for OFF (col B) and DEF (col C).
# ---------------------------------------------------------------------------
$ydsWs = $wb.Worksheets.Item("YDS")

$ydsWs.Range("B2").Value = $ydsWs.Range("B2").Value2 + " 6 3 4 6 6 -1 7 5 4 -5 5 -2 5 -4 1 2 3 11 4 0 2"
$ydsWs.Range("B3").Value = $ydsWs.Range("B3").Value2 + " 1 4 9 15 7 4 31 21 3 2 21 14 6 24 14 9 8 20 14 1 8 32 9 21 3 12 8 3 6 6 27 10 33"
$ydsWs.Range("C2").Value = $ydsWs.Range("C2").Value2 + " 55 2 0 1 -1 14 2 1 4 12 1 9 0 -1 0 1 0 0 11 7 2 4 22 -1 5"
$ydsWs.Range("C3").Value = $ydsWs.Range("C3").Value2 + " 9 9 13 5 9 12 5 9 24 11 24 19 2 23 11 21 14 4"

# ---------------------------------------------------------------------------
# OFF sheet: Home (row 2) / Road (row 3) season totals through Week 17.
# ---------------------------------------------------------------------------
$offWs = $wb.Worksheets.Item("OFF")

$offWs.Range("B2").Value = 6
$offWs.Range("C2").Value = 197
$offWs.Range("F2").Value = 25
$offWs.Range("G2").Value = 46
$offWs.Range("J2").Value = 31
$offWs.Range("N2").Value = 15
$offWs.Range("O2").Value = 15
$offWs.Range("P2").Value = 8

$offWs.Range("C3").Value = 207
$offWs.Range("E3").Value = 26
$offWs.Range("F3").Value = 158
$offWs.Range("G3").Value = 65
$offWs.Range("H3").Value = 18
$offWs.Range("I3").Value = 74
$offWs.Range("J3").Value = 64
$offWs.Range("L3").Value = 429
$offWs.Range("M3").Value = 290
$offWs.Range("Q3").Value = 666

# ---------------------------------------------------------------------------
# DEF sheet: Home (row 2) / Road (row 3) season totals through Week 17.
# ---------------------------------------------------------------------------
$defWs = $wb.Worksheets.Item("DEF")

$defWs.Range("B2").Value = 3
$defWs.Range("C2").Value = 162
$defWs.Range("D2").Value = 7
$defWs.Range("E2").Value = 8
$defWs.Range("F2").Value = 47
$defWs.Range("G2").Value = 46
$defWs.Range("J2").Value = 23
$defWs.Range("N2").Value = 31
$defWs.Range("O2").Value = 22

$defWs.Range("C3").Value = 229
$defWs.Range("D3").Value = 9
$defWs.Range("E3").Value = 34
$defWs.Range("F3").Value = 131
$defWs.Range("I3").Value = 67
$defWs.Range("J3").Value = 70
$defWs.Range("L3").Value = 367
$defWs.Range("M3").Value = 231
$defWs.Range("Q3").Value = 636

# ---------------------------------------------------------------------------
# ST sheet: season counters (row 2 / row 3) plus the per-kick/per-return
# logs (# column rows 4-6, RA column rows 3-5).
# ---------------------------------------------------------------------------
$stWs = $wb.Worksheets.Item("ST")

$stWs.Range("B2").Value = 95
$stWs.Range("D2").Value = 57
$stWs.Range("F2").Value = 459
$stWs.Range("G2").Value = 441
$stWs.Range("H2").Value = 3
$stWs.Range("I2").Value = 1
$stWs.Range("J2").Value = 187
$stWs.Range("K2").Value = 177
$stWs.Range("B3").Value = 73

$stWs.Range("B4").Value = $stWs.Range("B4").Value2 + " 63 66 61 60 63"
$stWs.Range("B5").Value = $stWs.Range("B5").Value2 + " 26 27 13 25 28"
$stWs.Range("B6").Value = $stWs.Range("B6").Value2 + " 7"
$stWs.Range("D3").Value = $stWs.Range("D3").Value2 + " 35 39 30"
$stWs.Range("D4").Value = $stWs.Range("D4").Value2 + " 0 0 0"
$stWs.Range("D5").Value = $stWs.Range("D5").Value2 + " 0 0 3 0"

# ---------------------------------------------------------------------------
# TURNS sheet: Road (row 3) turnover totals.
# ---------------------------------------------------------------------------
$turnsWs = $wb.Worksheets.Item("TURNS")

$turnsWs.Range("B3").Value = 9
$turnsWs.Range("E3").Value = 8

# ---------------------------------------------------------------------------
# PEN sheet: False start (row 2) penalty-yardage total.
# ---------------------------------------------------------------------------
$penWs = $wb.Worksheets.Item("PEN")

$penWs.Range("D2").Value = 17
